$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.dbnl.org/tekst/aren001kroo01_01"

$rows = @(
    @("#de-voorige", "De voorige"),
    @("#weeze", "Weeze"),
    @("#brittanje", "Brittanje"),
    @("#vleyery", "Vleyery"),
    @("#stroomnimf", "Stroomnimf"),
    @("#de", "De"),
    @("#teems", "Teems"),
    @("#seine", "Seine"),
    @("#de-maagd-van-groot-brittanje", "De Maagd van Groot Brittanje"),
    @("#weez", "Weez"),
    @("#e", "e"),
    @("#batavia", "Batavia"),
    @("#bedróg", "Bedróg"),
    @("#waarheid", "Waarheid"),
    @("#alle-te-zaamen", "Alle te zaamen"),
    @("#hovaardy", "Hovaardy")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $url
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 7).Value = ""
    $ws.Cells.Item($r, 8).Value = ""
    $r = $r + 1
}
